# Added API for BookCategory
# Replace the "Medical" branch (Medical / Surgery / Medicine / Gynae & Obs,
# rows 8-11) with five new "Engineering" children: Computer Science and
# Engineering2, Electrical Engineering2, Mechanical Engineering2, Chemical
# Engineering2, Civil Engineering2 (rows 8-12). Every row below shifts down
# by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Medical" category and its two direct children
# (rows 8:11 -> Medical / Surgery / Medicine / Gynae & Obs).
$ws.Range("A8:C11").EntireRow.Delete()

# Insert 5 fresh rows in their place for the new Engineering sub-categories.
$ws.Range("A8:C12").EntireRow.Insert()

# Match the row height used by every other data row on the sheet.
$ws.Range("A8:A12").EntireRow.RowHeight = 13.8

$ws.Cells.Item(8, 2).Value = "Computer Science and Engineering2"
$ws.Cells.Item(8, 3).Value = "Engineering"

$ws.Cells.Item(9, 2).Value = "Electrical Engineering2"
$ws.Cells.Item(9, 3).Value = "Engineering"

$ws.Cells.Item(10, 2).Value = "Mechanical Engineering2"
$ws.Cells.Item(10, 3).Value = "Engineering"

$ws.Cells.Item(11, 2).Value = "Chemical Engineering2"
$ws.Cells.Item(11, 3).Value = "Engineering"

$ws.Cells.Item(12, 2).Value = "Civil Engineering2"
$ws.Cells.Item(12, 3).Value = "Engineering"

# Match the author's final selection / scroll position.
$null = $ws.Range("A1").Select()
$null = $ws.Range("B12").Select()
